$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the zone_code value in row 2 (V2) from "CIV" to "COCO"
$ws.Range("V2").Value = "COCO"
